$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the three MAG rows that no longer appear in the output
# (delete bottom-to-top so row numbers of earlier rows stay stable)
$ws.Range("A10").EntireRow.Delete()
$ws.Range("A8").EntireRow.Delete()
$ws.Range("A6").EntireRow.Delete()

# Remove the "max" column (column C); D/E shift left into C/D
$ws.Range("C1").EntireColumn.Delete()

# New per-row prediction score / label values
$ws.Range("B2").Value = 0.01885178513548347
$ws.Range("C2").Value = "s__Clostridium_A leptum"
$ws.Range("D2").Value = "s__Clostridium_A leptum"

$ws.Range("B3").Value = -0.189169470233276
$ws.Range("C3").Value = "s__Clostridium_A leptum"
$ws.Range("D3").Value = "s__Clostridium_A leptum(reject)"

$ws.Range("B4").Value = -0.28751753218625
$ws.Range("C4").Value = "s__Clostridium_A leptum"
$ws.Range("D4").Value = "s__Clostridium_A leptum(reject)"

$ws.Range("B5").Value = -0.1570692448121331
$ws.Range("C5").Value = "s__Clostridium_A leptum"
$ws.Range("D5").Value = "s__Clostridium_A leptum(reject)"

$ws.Range("B6").Value = -0.1850295281991059
$ws.Range("C6").Value = "s__Clostridium_A leptum"
$ws.Range("D6").Value = "s__Clostridium_A leptum(reject)"

$ws.Range("B7").Value = -0.3097656684679784
$ws.Range("C7").Value = "s__Clostridium_A leptum"
$ws.Range("D7").Value = "s__Clostridium_A leptum(reject)"

$ws.Range("B8").Value = 0.04280660223353294
$ws.Range("C8").Value = "s__Clostridium_A leptum"
$ws.Range("D8").Value = "s__Clostridium_A leptum"
